$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update missing "cards" data -------------------------------------------------
# A handful of stories were saved with an empty placeholder cards list (`[""]`).
# Fill them in with the correct card lists, and add a new "unknow" entry where
# the real card could not be identified.
$ws.Range("E101").Value = '["Shokurei","Mishige"]'
$ws.Range("E102").Value = '["Shokurei","Mishige","Amezaiku"]'
$ws.Range("E107").Value = '["Inaba Kaguya","Mannendake","Ootengu","Hiromasa","unknow"]'
$ws.Range("E111").Value = '["Fukengaku","Momo","Sakura"]'
$ws.Range("E123").Value = '["Yamalord Enma","Hangan","Kuro Mujou","Shiro Mujou","Mouba"]'
$ws.Range("E124").Value = '["Heartseeker Momiji","Seimei"]'
$ws.Range("E135").Value = '["Ungaikyo","Zen Ungaikyo"]'

# --- Cosmetic touch-ups that came along with the content update -----------------
# Widen the "cards" column now that it holds longer lists, and make the
# "cards" header (E1) match the other (light-colored) table headers.
$ws.Columns("E").ColumnWidth = 69.375
$ws.Range("E2:E139").Font.ColorIndex = -4142
$ws.Range("E1").Font.ThemeColor = 2

# --- Restore the normal view (scroll position / selection) ----------------------
$ws.Range("H6").Select() | Out-Null
